$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tracker originally started on 2024-02-27 (serial 45349). Extend the
# cycle backward with one week of "period" (S) entries: 2024-02-20 through
# 2024-02-26 (serials 45342-45348), each flagged with phase status "S" in
# column C (Flujo), with no temperature reading for those days.
#
# We add these as brand-new rows below the current data (so we never touch
# the existing, untouched rows) and then sort the whole table by date so the
# new band takes its place at the top, pushing the existing readings down.

$newDates = 45342, 45343, 45344, 45345, 45346, 45347, 45348

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $startRow + $i
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $newDates[$i]
    $dateCell.NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
    $ws.Cells.Item($r, 3).Value = "S"
}

$lastRow2 = $startRow + $newDates.Length - 1
$fullRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow2, 3))
$keyRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow2, 1))
$fullRange.Sort($keyRange, 1)

# Restore the active selection to reflect where the new band now sits
[void]$ws.Range("D14").Select()

Write-Host "Inserted period band rows and sorted data by date."
